# Update simulation results on the "results" sheet (sheet1)
# after final checks / production run (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

$ws.Range("E2").Value = 9.081
$ws.Range("G2").Value = 11.897
$ws.Range("H2").Value = 1.792
$ws.Range("I2").Value = 1.102
$ws.Range("J2").Value = 0.05227754754146585
$ws.Range("K2").Value = 0.8155500711654082
$ws.Range("N2").Value = 1.705
$ws.Range("O2").Value = 1.056
$ws.Range("E3").Value = 9.025
$ws.Range("G3").Value = 11.86
$ws.Range("H3").Value = 1.159
$ws.Range("I3").Value = 0.742
$ws.Range("J3").Value = 0.048398393110158026
$ws.Range("K3").Value = 0.8136002271436689
$ws.Range("M3").Value = 110.065
$ws.Range("N3").Value = 1.22
$ws.Range("O3").Value = 0.767
$ws.Range("E4").Value = 8.977
$ws.Range("G4").Value = 11.651
$ws.Range("H4").Value = 1.587
$ws.Range("I4").Value = 1.055
$ws.Range("J4").Value = 0.05029629697488788
$ws.Range("K4").Value = 0.7866139846743296
$ws.Range("M4").Value = 114.441
$ws.Range("N4").Value = 1.647
$ws.Range("O4").Value = 0.984
$ws.Range("E5").Value = 9.184
$ws.Range("G5").Value = 11.977
$ws.Range("H5").Value = 2.821
$ws.Range("I5").Value = 1.727
$ws.Range("J5").Value = 0.04879827479850265
$ws.Range("K5").Value = 0.7997751481545005
$ws.Range("M5").Value = 123.538
$ws.Range("N5").Value = 2.749
$ws.Range("O5").Value = 1.704
$ws.Range("E6").Value = 8.926
$ws.Range("G6").Value = 15.037
$ws.Range("H6").Value = 2.048
$ws.Range("I6").Value = 1.154
$ws.Range("J6").Value = 0.07677565168285315
$ws.Range("K6").Value = 0.8880388218202115
$ws.Range("N6").Value = 2.158
$ws.Range("O6").Value = 1.139
$ws.Range("E7").Value = 9.019
$ws.Range("G7").Value = 14.804
$ws.Range("H7").Value = 1.565
$ws.Range("I7").Value = 0.809
$ws.Range("J7").Value = 0.07519247213055287
$ws.Range("K7").Value = 0.873874643874644
$ws.Range("M7").Value = 110.151
$ws.Range("N7").Value = 1.486
$ws.Range("E8").Value = 9.149
$ws.Range("G8").Value = 15.06
$ws.Range("H8").Value = 2.116
$ws.Range("I8").Value = 1.135
$ws.Range("J8").Value = 0.07783003924957392
$ws.Range("K8").Value = 0.8852601531075642
$ws.Range("M8").Value = 114.496
$ws.Range("N8").Value = 2.223
$ws.Range("O8").Value = 1.189
$ws.Range("E9").Value = 8.973
$ws.Range("G9").Value = 14.961
$ws.Range("H9").Value = 3.685
$ws.Range("I9").Value = 1.829
$ws.Range("J9").Value = 0.07819802477430075
$ws.Range("K9").Value = 0.8850196142530234
$ws.Range("M9").Value = 123.673
$ws.Range("N9").Value = 3.602
$ws.Range("O9").Value = 1.931
$ws.Range("E10").Value = 8.865
$ws.Range("G10").Value = 18.116
$ws.Range("H10").Value = 2.549
$ws.Range("I10").Value = 1.27
$ws.Range("J10").Value = 0.1043663122383385
$ws.Range("K10").Value = 0.964919141424668
$ws.Range("N10").Value = 2.751
$ws.Range("O10").Value = 1.263
$ws.Range("E11").Value = 8.839
$ws.Range("G11").Value = 18.027
$ws.Range("H11").Value = 1.818
$ws.Range("I11").Value = 0.861
$ws.Range("J11").Value = 0.10024594351297163
$ws.Range("K11").Value = 0.957715540911587
$ws.Range("M11").Value = 110.175
$ws.Range("N11").Value = 1.83
$ws.Range("O11").Value = 0.899
$ws.Range("E12").Value = 9.033
$ws.Range("G12").Value = 18.288
$ws.Range("H12").Value = 2.491
$ws.Range("I12").Value = 1.242
$ws.Range("J12").Value = 0.10511188200154083
$ws.Range("K12").Value = 0.9681407513076559
$ws.Range("M12").Value = 114.339
$ws.Range("N12").Value = 2.573
$ws.Range("O12").Value = 1.207
$ws.Range("E13").Value = 9.216
$ws.Range("G13").Value = 18.448
$ws.Range("H13").Value = 4.427
$ws.Range("I13").Value = 1.941
$ws.Range("J13").Value = 0.10689618361214619
$ws.Range("K13").Value = 0.9705815958335553
$ws.Range("M13").Value = 123.454
$ws.Range("N13").Value = 4.426
$ws.Range("O13").Value = 2.166
$ws.Range("E14").Value = 9.097
$ws.Range("G14").Value = 19.691
$ws.Range("H14").Value = 3.003
$ws.Range("I14").Value = 1.224
$ws.Range("J14").Value = 0.1181841012767193
$ws.Range("K14").Value = 0.9904959978749329
$ws.Range("N14").Value = 3.1
$ws.Range("O14").Value = 1.331
$ws.Range("E15").Value = 9.033
$ws.Range("G15").Value = 19.563
$ws.Range("H15").Value = 2.013
$ws.Range("I15").Value = 0.91
$ws.Range("J15").Value = 0.12434547837353586
$ws.Range("K15").Value = 0.9818965517241379
$ws.Range("M15").Value = 110.151
$ws.Range("N15").Value = 2.008
$ws.Range("O15").Value = 0.859
$ws.Range("E16").Value = 9.035
$ws.Range("G16").Value = 19.708
$ws.Range("H16").Value = 2.967
$ws.Range("I16").Value = 1.287
$ws.Range("J16").Value = 0.11268901882635536
$ws.Range("K16").Value = 0.9879629629629629
$ws.Range("M16").Value = 114.345
$ws.Range("N16").Value = 2.755
$ws.Range("O16").Value = 1.269
$ws.Range("E17").Value = 8.993
$ws.Range("G17").Value = 19.498
$ws.Range("H17").Value = 4.303
$ws.Range("I17").Value = 2.144
$ws.Range("J17").Value = 0.1153397263243248
$ws.Range("K17").Value = 0.990478994895956
$ws.Range("M17").Value = 123.616
$ws.Range("N17").Value = 4.59
$ws.Range("O17").Value = 2.086
